$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 / rId1) - update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 918
$ws1.Range("F7").Value = 4505
$ws1.Range("F8").Value = 2672
$ws1.Range("F10").Value = 2633
$ws1.Range("F14").Value = 1682
$ws1.Range("F16").Value = 201
$ws1.Range("F18").Value = 354
$ws1.Range("F19").Value = 33
$ws1.Range("F22").Value = 49
$ws1.Range("F23").Value = 496
$ws1.Range("F24").Value = 35
$ws1.Range("F26").Value = 591
$ws1.Range("F27").Value = 719
$ws1.Range("F28").Value = 128
$ws1.Range("F29").Value = 14
$ws1.Range("F30").Value = 458
$ws1.Range("F31").Value = 1637
$ws1.Range("F32").Value = 1250
$ws1.Range("F33").Value = 233
$ws1.Range("F34").Value = 32
$ws1.Range("F35").Value = 1307
$ws1.Range("F36").Value = 2168
$ws1.Range("F37").Value = 328
$ws1.Range("F38").Value = 18
$ws1.Range("F39").Value = 569
$ws1.Range("F40").Value = 101
$ws1.Range("F43").Value = 715
$ws1.Range("F44").Value = 1391
$ws1.Range("F45").Value = 152
$ws1.Range("F47").Value = 456
$ws1.Range("F49").Value = 87

# Sheet "演出" (sheet2 / rId2) - update column F (想去人数) values
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F14").Value = 13

# Sheet "全部类型" (sheet4 / rId4) - update column F (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 918
$ws4.Range("F5").Value = 4505
$ws4.Range("F6").Value = 2672
$ws4.Range("F7").Value = 2633
$ws4.Range("F8").Value = 1682
$ws4.Range("F12").Value = 201
$ws4.Range("F14").Value = 354
$ws4.Range("F15").Value = 33
$ws4.Range("F18").Value = 49
$ws4.Range("F19").Value = 496
$ws4.Range("F20").Value = 35
$ws4.Range("F21").Value = 591
$ws4.Range("F22").Value = 719
$ws4.Range("F23").Value = 128
$ws4.Range("F27").Value = 458
$ws4.Range("F28").Value = 1637
$ws4.Range("F29").Value = 1250
$ws4.Range("F30").Value = 233
$ws4.Range("F31").Value = 32
$ws4.Range("F34").Value = 2168
$ws4.Range("F35").Value = 328
$ws4.Range("F39").Value = 569
$ws4.Range("F40").Value = 101
$ws4.Range("F43").Value = 715
$ws4.Range("F44").Value = 1391
$ws4.Range("F46").Value = 152
$ws4.Range("F47").Value = 456
$ws4.Range("F48").Value = 87
$ws4.Range("F49").Value = 13
